$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.3316315710544586
$ws.Range("B1").Value = 2.403627157211304
$ws.Range("C1").Value = 4.781197547912598
$ws.Range("D1").Value = 1.687149286270142
$ws.Range("E1").Value = 0.8499595522880554
